$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new worksheet "iterative fungi search" as the last sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "iterative fungi search"

# ---------------------------------------------------------------------------
# Row 3: header row for first table (entered first)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Sample"
$ws.Range("B3").Value = "0 MC peptides"
$ws.Range("C3").Value = "1 MC peptides"
$ws.Range("D3").Value = "2 MC peptides"

# ---------------------------------------------------------------------------
# Row 1: section title (entered after the header row)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Unique fungal peptides found by iterative database searching:"

# ---------------------------------------------------------------------------
# Second table: header label + first new-text column, then first block of rows
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Sample"
$ws.Range("B12").Value = "Depth (m)"

$ws.Range("A13").Value = "Epipelagic suspended"
$ws.Range("A14").Value = "Mesopelagic suspended"
$ws.Range("A15").Value = "Bathypelagic suspended"

# ---------------------------------------------------------------------------
# Second table: the "missed cleavages" column header, then second block of rows
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "Missed tryptic cleaveages:"

$ws.Range("A16").Value = "Epipelagic sinking"
$ws.Range("A17").Value = "Mesopelagic sinking"
$ws.Range("A18").Value = "Bathypelagic sinking"

# ---------------------------------------------------------------------------
# Rows 4-9: first table numeric data
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 231
$ws.Range("B4").Value = 17
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0

$ws.Range("A5").Value = 233
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 0

$ws.Range("A6").Value = 243
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 0

$ws.Range("A7").Value = 378
$ws.Range("B7").Value = 18
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 0

$ws.Range("A8").Value = 278
$ws.Range("B8").Value = 45
$ws.Range("C8").Value = 9
$ws.Range("D8").Value = 1

$ws.Range("A9").Value = 273
$ws.Range("B9").Value = 34
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 0

# ---------------------------------------------------------------------------
# Row 12: numeric column headers (0/1/2) + bottom border style
# ---------------------------------------------------------------------------
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("A12:F12").Borders.Item(9).LineStyle = 1
$ws.Range("A12:F12").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# Rows 13-18: second table numeric data
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = 100
$ws.Range("D13").Value = 17
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0

$ws.Range("B14").Value = 265
$ws.Range("D14").Value = 25
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 0

$ws.Range("B15").Value = 1000
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0

$ws.Range("B16").Value = 100
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 0

$ws.Range("B17").Value = 265
$ws.Range("D17").Value = 45
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 1

$ws.Range("B18").Value = 965
$ws.Range("D18").Value = 34
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 0

# ---------------------------------------------------------------------------
# Column widths (closest achievable values matching target OOXML widths of
# 21.6640625 / 22.33203125 characters)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.833333333333336
$ws.Columns.Item(3).ColumnWidth = 21.5

# ---------------------------------------------------------------------------
# Page setup (portrait orientation, matching the target sheet)
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selection / view state for the new sheet
# ---------------------------------------------------------------------------
$ws.Range("G16").Select()

# ---------------------------------------------------------------------------
# Update the previously-active sheet ("fungi dno + db peps tryptic") view:
# it is no longer the tab-selected sheet, and its scroll/selection changes.
# ---------------------------------------------------------------------------
$fungiSheet = $wb.Worksheets.Item("fungi dno + db peps tryptic")
$fungiSheet.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$fungiSheet.Range("J17").Select()

# ---------------------------------------------------------------------------
# Make the new sheet the active tab, matching the target workbook view
# ---------------------------------------------------------------------------
$ws.Activate()
